$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates: rows 2-3 move to the later date, rows 4-5 move to the earlier date
$ws.Range("D2").Value = 44574
$ws.Range("D3").Value = 44574
$ws.Range("D4").Value = 44559
$ws.Range("D5").Value = 44559
